$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 12500581
$ws.Range("I107").Value = 16666988
$ws.Range("K107").Value = 16666988
$ws.Range("M107").Value = -16665068

$ws.Range("H116").Value = 11225.417
$ws.Range("I116").Value = 21601
$ws.Range("J116").Value = 3814.2856
$ws.Range("K116").Value = 21601
$ws.Range("L116").Value = 3814.2856
$ws.Range("M116").Value = -18159
$ws.Range("N116").Value = -10698.2856

$ws.Range("H129").Value = 874.3514
$ws.Range("I129").Value = 488.29413
$ws.Range("J129").Value = 989.4912
$ws.Range("K129").Value = 1464.88239
$ws.Range("L129").Value = 2968.4736
$ws.Range("M129").Value = 3535.11761
$ws.Range("N129").Value = -12968.4736

$ws.Range("H137").Value = 1424.5588
$ws.Range("I137").Value = 1292.4073
$ws.Range("J137").Value = 1934.2858
$ws.Range("K137").Value = 3877.2219
$ws.Range("L137").Value = 5802.857400000001
$ws.Range("M137").Value = -1327.2219
$ws.Range("N137").Value = -10902.8574

$ws.Range("H138").Value = 2710.04
$ws.Range("I138").Value = 850.5
$ws.Range("J138").Value = 3849.758
$ws.Range("K138").Value = 2551.5
$ws.Range("L138").Value = 11549.274
$ws.Range("M138").Value = 2588.5
$ws.Range("N138").Value = -21829.274

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4723.1665
$ws.Range("I32").Value = 2929.4666
$ws.Range("K32").Value = 2929.4666
$ws.Range("M32").Value = -2642.4666

$ws.Range("H61").Value = 3314.2307
$ws.Range("I61").Value = 3540.111
$ws.Range("K61").Value = 3540.111
$ws.Range("M61").Value = -3328.111

$ws.Range("H74").Value = 1873.0416
$ws.Range("I74").Value = 1664.9375
$ws.Range("J74").Value = 2289.25
$ws.Range("K74").Value = 1664.9375
$ws.Range("L74").Value = 2289.25
$ws.Range("M74").Value = -790.9375
$ws.Range("N74").Value = -4037.25

$ws.Range("H77").Value = 1873.0416
$ws.Range("I77").Value = 1664.9375
$ws.Range("J77").Value = 2289.25
$ws.Range("K77").Value = 8324.6875
$ws.Range("L77").Value = 11446.25
$ws.Range("M77").Value = -3956.6875
$ws.Range("N77").Value = -20182.25

$ws.Range("H132").Value = 2696.6592
$ws.Range("I132").Value = 1705.6428
$ws.Range("J132").Value = 4430.9375
$ws.Range("K132").Value = 5116.928400000001
$ws.Range("L132").Value = 13292.8125
$ws.Range("M132").Value = -2586.928400000001
$ws.Range("N132").Value = -18352.8125

$ws.Range("H136").Value = 3314.2307
$ws.Range("I136").Value = 3540.111
$ws.Range("K136").Value = 10620.333
$ws.Range("M136").Value = -8070.332999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 33335144
$ws.Range("I86").Value = 41668416
$ws.Range("J86").Value = 2050
$ws.Range("K86").Value = 41668416
$ws.Range("L86").Value = 2050
$ws.Range("M86").Value = -41667293
$ws.Range("N86").Value = -4296

$ws.Range("H89").Value = 33335144
$ws.Range("I89").Value = 41668416
$ws.Range("J89").Value = 2050
$ws.Range("K89").Value = 208342080
$ws.Range("L89").Value = 10250
$ws.Range("M89").Value = -208336464
$ws.Range("N89").Value = -21482

$ws.Range("H94").Value = 2096.36
$ws.Range("I94").Value = 1806.8125
$ws.Range("J94").Value = 2611.111
$ws.Range("K94").Value = 1806.8125
$ws.Range("L94").Value = 2611.111
$ws.Range("M94").Value = -1355.8125
$ws.Range("N94").Value = -3513.111

$ws.Range("H99").Value = 166667820
$ws.Range("I99").Value = 200001060
$ws.Range("K99").Value = 200001060
$ws.Range("M99").Value = -199999562

$ws.Range("H134").Value = 6119.6294
$ws.Range("I134").Value = 9475.071
$ws.Range("J134").Value = 2506.077
$ws.Range("K134").Value = 28425.213
$ws.Range("L134").Value = 7518.231000000001
$ws.Range("M134").Value = -25890.213
$ws.Range("N134").Value = -12588.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9266.666999999999
$ws.Range("J4").Value = 7148.148
$ws.Range("L4").Value = 7148.148
$ws.Range("N4").Value = -7372.148

$ws.Range("H31").Value = 3071.889
$ws.Range("I31").Value = 1457.8684
$ws.Range("J31").Value = 5525.2
$ws.Range("K31").Value = 1457.8684
$ws.Range("L31").Value = 5525.2
$ws.Range("M31").Value = -1162.8684
$ws.Range("N31").Value = -6115.2

$ws.Range("H34").Value = 3071.889
$ws.Range("I34").Value = 1457.8684
$ws.Range("J34").Value = 5525.2
$ws.Range("K34").Value = 1457.8684
$ws.Range("L34").Value = 5525.2
$ws.Range("M34").Value = -1255.8684
$ws.Range("N34").Value = -5929.2

$ws.Range("H122").Value = 2246.0454
$ws.Range("I122").Value = 1795.2667
$ws.Range("K122").Value = 5385.800099999999
$ws.Range("M122").Value = -2935.800099999999

$ws.Range("H134").Value = 2026.8959
$ws.Range("I134").Value = 2230.7297
$ws.Range("K134").Value = 6692.1891
$ws.Range("M134").Value = -4157.1891

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1821.4546
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 1821.4546
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 5464.3638
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -7086.3638

$ws.Range("H72").Value = 1821.4546
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 1821.4546
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 16393.0914
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -24505.0914

$ws.Range("H100").Value = 3483.3333
$ws.Range("I100").Value = 3500
$ws.Range("J100").Value = 3480
$ws.Range("K100").Value = 10500
$ws.Range("L100").Value = 10440
$ws.Range("M100").Value = -9689
$ws.Range("N100").Value = -12062

$ws.Range("H103").Value = 2839.6
$ws.Range("I103").Value = 999.5
$ws.Range("J103").Value = 3299.625
$ws.Range("K103").Value = 2998.5
$ws.Range("L103").Value = 9898.875
$ws.Range("M103").Value = -2119.5
$ws.Range("N103").Value = -11656.875

$ws.Range("H125").Value = 3588.923
$ws.Range("I125").Value = 993
$ws.Range("J125").Value = 4367.7
$ws.Range("K125").Value = 2979
$ws.Range("L125").Value = 13103.1
$ws.Range("M125").Value = 1941
$ws.Range("N125").Value = -22943.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 22807.047
$ws.Range("J123").Value = 22807.047
$ws.Range("L123").Value = 22807.047
$ws.Range("N123").Value = -27707.047

$ws.Range("H132").Value = 2284.0732
$ws.Range("I132").Value = 2132.2222
$ws.Range("J132").Value = 2576.9285
$ws.Range("K132").Value = 6396.6666
$ws.Range("L132").Value = 7730.7855
$ws.Range("M132").Value = -3866.6666
$ws.Range("N132").Value = -12790.7855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2199.6667
$ws.Range("J4").Value = 2199.6667
$ws.Range("L4").Value = 2199.6667
$ws.Range("N4").Value = -2425.6667

$ws.Range("H15").Value = 5500
$ws.Range("J15").Value = 5500
$ws.Range("L15").Value = 5500
$ws.Range("N15").Value = -6076

$ws.Range("H96").Value = 2747.8
$ws.Range("I96").Value = 2529
$ws.Range("J96").Value = 3258.3333
$ws.Range("K96").Value = 2529
$ws.Range("L96").Value = 3258.3333
$ws.Range("M96").Value = -1156
$ws.Range("N96").Value = -6004.3333
